$wb = $excel.ActiveWorkbook

# ---- Sheet: P_valores ----
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.8154008458012483
$wsP.Range("D2").Value = 0.5347106467231812
$wsP.Range("E2").Value = 0.8184909588799254
$wsP.Range("F2").Value = 0.7321273681438538

$wsP.Range("B3").Value = 0.8154008458012483
$wsP.Range("D3").Value = 0.4560039731195333
$wsP.Range("E3").Value = 0.5293652828311641
$wsP.Range("F3").Value = 0.5757957273994081

$wsP.Range("B4").Value = 0.5347106467231812
$wsP.Range("C4").Value = 0.4560039731195333
$wsP.Range("E4").Value = 0.6450878918800242
$wsP.Range("F4").Value = 0.5975588300474062

$wsP.Range("B5").Value = 0.8184909588799254
$wsP.Range("C5").Value = 0.5293652828311641
$wsP.Range("D5").Value = 0.6450878918800242
$wsP.Range("F5").Value = 0.9115069976567494

$wsP.Range("B6").Value = 0.7321273681438538
$wsP.Range("C6").Value = 0.5757957273994081
$wsP.Range("D6").Value = 0.5975588300474062
$wsP.Range("E6").Value = 0.9115069976567494

# ---- Sheet: Estadisticos_DM ----
$wsD = $wb.Worksheets.Item("Estadisticos_DM")

$wsD.Range("C2").Value = -0.2379056700193651
$wsD.Range("D2").Value = 0.6365141036749664
$wsD.Range("E2").Value = 0.2338440244291884
$wsD.Range("F2").Value = 0.3492113672464949

$wsD.Range("B3").Value = 0.2379056700193651
$wsD.Range("D3").Value = 0.7666840612921585
$wsD.Range("E3").Value = 0.6449850667625696
$wsD.Range("F3").Value = 0.5729047729763265

$wsD.Range("B4").Value = -0.6365141036749664
$wsD.Range("C4").Value = -0.7666840612921585
$wsD.Range("E4").Value = -0.470724001009268
$wsD.Range("F4").Value = -0.5401842302751506

$wsD.Range("B5").Value = -0.2338440244291884
$wsD.Range("C5").Value = -0.6449850667625696
$wsD.Range("D5").Value = 0.470724001009268
$wsD.Range("F5").Value = 0.1131643266339309

$wsD.Range("B6").Value = -0.3492113672464949
$wsD.Range("C6").Value = -0.5729047729763265
$wsD.Range("D6").Value = 0.5401842302751506
$wsD.Range("E6").Value = -0.1131643266339309

$wb.Save()
